$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, shifting existing rows 31-146 down to 32-147
$ws.Rows.Item(31).Insert()

# Populate the new row 31 with the new weekly record
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 44910
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100108
$ws.Range("H31").Value = "Tropicales y subtropicales"
$ws.Range("I31").Value = 100108002
$ws.Range("J31").Value = "Mango"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 7500
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 7667
$ws.Range("Q31").Value = "$/bandeja 4 kilos"
$ws.Range("R31").Value = "Perú"
$ws.Range("S31").Value = 1917
$ws.Range("T31").Value = 4
